$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix schema/column-name typos: the camelCase DB column names were
# replaced with their actual (lowercase) MySQL column names.
# Order matters for shared-string table layout (matches how Excel
# appends newly-introduced strings as it re-saves the sheet).
$ws.Range("A13").Value = "namecd"
$ws.Range("A15").Value = "nametype"
$ws.Range("A14").Value = "nameeng"
$ws.Range("A10").Value = "deletereason"

# Autofit column A now that its content changed, and leave the cursor
# on D11 (matches the script's last selection).
$ws.Columns("A:A").AutoFit()
$null = $ws.Range("D11").Select()
